$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column U ("Trans") to "N" for the data rows that were missing a value:
# rows 2-45 and rows 77-140 (rows 46-76 already contain "Y" and are left untouched)
$ws.Range("U2:U45").Value = "N"
$ws.Range("U77:U140").Value = "N"
